# Actualización automática 2025-06-30 14:35:09
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (per-client, per-product-group sales)
# ------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M3").Value  = 1317.8               # BECERRA FARIAS ROSA DAYANA - PORCELANATO
$wsGrupo.Range("M10").Value = 11.52                # F.V - AREA ANDINA S.A. - PORCELANATO
$wsGrupo.Range("E15").Value = 64.81999999999999    # TOSCANO RAMIREZ MONICA CECILIA - FREGADEROS DE COCINA
$wsGrupo.Range("G15").Value = 40.74                # TOSCANO RAMIREZ MONICA CECILIA - GRIFERIAS

# Row 19 totals ("N de 17" counters)
$wsGrupo.Range("E19").Value = "1 de 17"
$wsGrupo.Range("G19").Value = "1 de 17"
$wsGrupo.Range("M19").Value = "4 de 17"

# ------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (monthly sales, junio column = F)
# ------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F3").Value  = 1317.8             # BECERRA FARIAS ROSA DAYANA
$wsMensual.Range("F10").Value = 11.52              # F.V - AREA ANDINA S.A.
$wsMensual.Range("F15").Value = 2206.26            # TOSCANO RAMIREZ MONICA CECILIA
$wsMensual.Range("F19").Value = 33991.41           # total

# ------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (budget vs. sale vs. compliance, per product group)
# ------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# FREGADEROS DE COCINA (row 4)
$wsCumpl.Range("D4").Value = 64.81999999999999
$wsCumpl.Range("E4").Value = 308.173863046034
$wsCumpl.Range("F4").Value = 0.1737830201029341

# GRIFERIAS (row 6)
$wsCumpl.Range("D6").Value = 40.74
$wsCumpl.Range("E6").Value = 66.07999999999998
$wsCumpl.Range("F6").Value = 0.381389252948886

# PORCELANATO (row 16)
$wsCumpl.Range("D16").Value = 30989.09
$wsCumpl.Range("E16").Value = -2779.25
$wsCumpl.Range("F16").Value = 1.098520587142642

# TOTAL (row 19)
$wsCumpl.Range("D19").Value = 33991.41
$wsCumpl.Range("E19").Value = 13227.89386304603
$wsCumpl.Range("F19").Value = 0.7198625820191682
